$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1599
$ws1.Range("F4").Value = 8540
$ws1.Range("F6").Value = 63
$ws1.Range("F8").Value = 81
$ws1.Range("F9").Value = 1333
$ws1.Range("F10").Value = 107
$ws1.Range("F13").Value = 9188
$ws1.Range("F19").Value = 6085
$ws1.Range("F20").Value = 1043
$ws1.Range("F21").Value = 60
$ws1.Range("F23").Value = 103

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1599
$ws4.Range("F4").Value = 8540
$ws4.Range("F6").Value = 63
$ws4.Range("F8").Value = 81
$ws4.Range("F9").Value = 1333
$ws4.Range("F10").Value = 107
$ws4.Range("F15").Value = 9188
$ws4.Range("F21").Value = 6085
$ws4.Range("F22").Value = 1043
$ws4.Range("F23").Value = 60
$ws4.Range("F25").Value = 103
